# Update cryptos list data
# Commit message: Updated cryptos list on Tue Feb 20 05:42:49 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '51.976.85'
$ws.Cells.Item(2, 5).Value = '  -0.13%  '
$ws.Cells.Item(3, 4).Value = '2.928.78'
$ws.Cells.Item(3, 5).Value = '  +1.40%  '
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$ws.Cells.Item(5, 4).Value = '''357.13'
$ws.Cells.Item(5, 5).Value = '  +1.70%  '
$ws.Cells.Item(6, 4).Value = '''111.01'
$ws.Cells.Item(6, 5).Value = '  -0.32%  '
$ws.Cells.Item(7, 4).Value = '''0.568'
$ws.Cells.Item(7, 5).Value = '  +1.68%  '
$ws.Cells.Item(8, 5).Value = '  -0.06%  '
$ws.Cells.Item(9, 4).Value = '''0.629'
$ws.Cells.Item(9, 5).Value = '  +1.19%  '
$ws.Cells.Item(10, 4).Value = '''39.50'
$ws.Cells.Item(10, 5).Value = '  -0.81%  '
$ws.Cells.Item(11, 4).Value = '''0.0881'
$ws.Cells.Item(11, 5).Value = '  +2.93%  '
$ws.Cells.Item(12, 4).Value = '''0.137'
$ws.Cells.Item(12, 5).Value = '  +0.61%  '
$ws.Cells.Item(13, 4).Value = '''19.71'
$ws.Cells.Item(13, 5).Value = '  -1.06%  '
$ws.Cells.Item(14, 4).Value = '''7.91'
$ws.Cells.Item(14, 5).Value = '  +1.97%  '
$ws.Cells.Item(15, 4).Value = '3.392.38'
$ws.Cells.Item(15, 5).Value = '  +1.27%  '
$ws.Cells.Item(16, 4).Value = '2.918.84'
$ws.Cells.Item(16, 5).Value = '  +1.08%  '
$ws.Cells.Item(17, 4).Value = '''0.988'
$ws.Cells.Item(17, 5).Value = '  -1.02%  '
$ws.Cells.Item(18, 4).Value = '51.963.27'
$ws.Cells.Item(18, 5).Value = '  -0.26%  '
$ws.Cells.Item(19, 4).Value = '''3.31'
$ws.Cells.Item(19, 5).Value = '  -0.31%  '
$ws.Cells.Item(20, 4).Value = '''7.58'
$ws.Cells.Item(20, 5).Value = '  -1.47%  '
$ws.Cells.Item(21, 4).Value = '''14.06'
$ws.Cells.Item(21, 5).Value = '  -3.01%  '
$ws.Cells.Item(22, 4).Value = '0.0₃0983'
$ws.Cells.Item(22, 5).Value = '  +0.45%  '
$ws.Cells.Item(23, 4).Value = '''71.03'
$ws.Cells.Item(23, 5).Value = '  +0.49%  '
$ws.Cells.Item(24, 4).Value = '''270.80'
$ws.Cells.Item(24, 5).Value = '  +0.54%  '
$ws.Cells.Item(25, 4).Value = '''2.82'
$ws.Cells.Item(25, 5).Value = '  +1.55%  '
$ws.Cells.Item(26, 5).Value = '  +12.51%  '
$ws.Cells.Item(27, 4).Value = '''27.21'
$ws.Cells.Item(27, 5).Value = '  +2.72%  '
$ws.Cells.Item(28, 5).Value = '  +0.17%  '
$ws.Cells.Item(29, 5).Value = '  +14.91%  '
$ws.Cells.Item(30, 4).Value = '''0.107'
$ws.Cells.Item(30, 5).Value = '  +13.27%  '
$ws.Cells.Item(31, 5).Value = '  +0.95%  '
$ws.Cells.Item(32, 4).Value = '''38.79'
$ws.Cells.Item(32, 5).Value = '  +1.12%  '
$ws.Cells.Item(33, 4).Value = '''6.05'
$ws.Cells.Item(33, 5).Value = '  -1.66%  '
$ws.Cells.Item(34, 4).Value = '''52.20'
$ws.Cells.Item(34, 5).Value = '  -1.24%  '
$ws.Cells.Item(35, 4).Value = '''0.0446'
$ws.Cells.Item(35, 5).Value = '  -2.44%  '
$ws.Cells.Item(36, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(36, 4).Value = '''0.999'
$ws.Cells.Item(36, 5).Value = '  -0.05%  '
$ws.Cells.Item(37, 2).Value = 'Toncoin'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(37, 4).Value = '''1.92'
$ws.Cells.Item(37, 5).Value = '  -14.05%  '
$ws.Cells.Item(38, 4).Value = '''3.24'
$ws.Cells.Item(38, 5).Value = '  -1.49%  '
$ws.Cells.Item(39, 5).Value = '  -0.07%  '
$ws.Cells.Item(40, 5).Value = '  -0.84%  '
$ws.Cells.Item(41, 4).Value = '''2.76'
$ws.Cells.Item(41, 5).Value = '  +4.60%  '
$ws.Cells.Item(42, 5).Value = '  +3.16%  '
$ws.Cells.Item(43, 4).Value = '''23.24'
$ws.Cells.Item(43, 5).Value = '  +2.36%  '
$ws.Cells.Item(44, 4).Value = '''119.35'
$ws.Cells.Item(44, 5).Value = '  -2.25%  '
$ws.Cells.Item(45, 5).Value = '  -1.44%  '
$ws.Cells.Item(46, 5).Value = '  +0.26%  '
$ws.Cells.Item(47, 5).Value = '  -2.72%  '
$ws.Cells.Item(48, 4).Value = '2.142.37'
$ws.Cells.Item(48, 5).Value = '  -2.51%  '
$ws.Cells.Item(49, 5).Value = '  -8.61%  '
$ws.Cells.Item(50, 4).Value = '''0.0333'
$ws.Cells.Item(50, 5).Value = '  +2.77%  '
$ws.Cells.Item(51, 2).Value = 'FraxShare'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(51, 4).Value = '''9.22'
$ws.Cells.Item(51, 5).Value = '  +1.58%  '
